{"js": "async (context) => {\n  // Map of old text -> new text, in document order.\n  // Index 0 is the date heading; 1..25 are the multiplication problems\n  // inside the table cells.\n  const replacements = [\n    [\"2025-11-26 Wednesday\", \"2025-11-27 Thursday\"],\n    [\"334\u00d72=668\", \"530\u00d76=3180\"],\n    [\"563\u00d79=5067\", \"853\u00d77=5971\"],\n    [\"566\u00d74=2264\", \"373\u00d79=3357\"],\n    [\"947\u00d79=8523\", \"563\u00d72=1126\"],\n    [\"576\u00d77=4032\", \"636\u00d73=1908\"],\n    [\"481\u00d72=962\", \"872\u00d77=6104\"],\n    [\"389\u00d72=778\", \"477\u00d79=4293\"],\n    [\"405\u00d78=3240\", \"902\u00d72=1804\"],\n    [\"272\u00d72=544\", \"221\u00d75=1105\"],\n    [\"820\u00d75=4100\", \"864\u00d77=6048\"],\n    [\"231\u00d77=1617\", \"794\u00d75=3970\"],\n    [\"298\u00d77=2086\", \"635\u00d74=2540\"],\n    [\"135\u00d77=945\", \"718\u00d74=2872\"],\n    [\"973\u00d72=1946\", \"251\u00d79=2259\"],\n    [\"950\u00d77=6650\", \"730\u00d77=5110\"],\n    [\"709\u00d78=5672\", \"239\u00d76=1434\"],\n    [\"846\u00d78=6768\", \"789\u00d73=2367\"],\n    [\"274\u00d75=1370\", \"763\u00d75=3815\"],\n    [\"259\u00d72=518\", \"617\u00d75=3085\"],\n    [\"957\u00d77=6699\", \"181\u00d74=724\"],\n    [\"347\u00d73=1041\", \"913\u00d77=6391\"],\n    [\"911\u00d72=1822\", \"417\u00d75=2085\"],\n    [\"788\u00d76=4728\", \"737\u00d75=3685\"],\n    [\"707\u00d77=4949\", \"337\u00d72=674\"],\n    [\"715\u00d79=6435\", \"495\u00d78=3960\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      throw new Error(`Text not found: ${oldText}`);\n    }\n\n    for (const range of results.items) {\n      range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old text -> new text, in document order.\n# Index 0 is the date heading; the rest are the multiplication\n# problems inside the table cells. All source strings are unique\n# within the document, so a straightforward Find/ReplaceAll per pair\n# is unambiguous.\n$replacements = @(\n    @(\"2025-11-26 Wednesday\", \"2025-11-27 Thursday\"),\n    @(\"334\u00d72=668\", \"530\u00d76=3180\"),\n    @(\"563\u00d79=5067\", \"853\u00d77=5971\"),\n    @(\"566\u00d74=2264\", \"373\u00d79=3357\"),\n    @(\"947\u00d79=8523\", \"563\u00d72=1126\"),\n    @(\"576\u00d77=4032\", \"636\u00d73=1908\"),\n    @(\"481\u00d72=962\", \"872\u00d77=6104\"),\n    @(\"389\u00d72=778\", \"477\u00d79=4293\"),\n    @(\"405\u00d78=3240\", \"902\u00d72=1804\"),\n    @(\"272\u00d72=544\", \"221\u00d75=1105\"),\n    @(\"820\u00d75=4100\", \"864\u00d77=6048\"),\n    @(\"231\u00d77=1617\", \"794\u00d75=3970\"),\n    @(\"298\u00d77=2086\", \"635\u00d74=2540\"),\n    @(\"135\u00d77=945\", \"718\u00d74=2872\"),\n    @(\"973\u00d72=1946\", \"251\u00d79=2259\"),\n    @(\"950\u00d77=6650\", \"730\u00d77=5110\"),\n    @(\"709\u00d78=5672\", \"239\u00d76=1434\"),\n    @(\"846\u00d78=6768\", \"789\u00d73=2367\"),\n    @(\"274\u00d75=1370\", \"763\u00d75=3815\"),\n    @(\"259\u00d72=518\", \"617\u00d75=3085\"),\n    @(\"957\u00d77=6699\", \"181\u00d74=724\"),\n    @(\"347\u00d73=1041\", \"913\u00d77=6391\"),\n    @(\"911\u00d72=1822\", \"417\u00d75=2085\"),\n    @(\"788\u00d76=4728\", \"737\u00d75=3685\"),\n    @(\"707\u00d77=4949\", \"337\u00d72=674\"),\n    @(\"715\u00d79=6435\", \"495\u00d78=3960\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n$d.Save()\n"}
